$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.100.84"
$ws.Range("E2").Value = "'  -0.84%  "
$ws.Range("D3").Value = "'1.909.58"
$ws.Range("E3").Value = "'  -1.31%  "
$ws.Range("E4").Value = "'  -0.22%  "
$ws.Range("D5").Value = "'0.7422"
$ws.Range("E5").Value = "'  -0.98%  "
$ws.Range("D6").Value = "'244.61"
$ws.Range("E6").Value = "'  +0.36%  "
$ws.Range("E7").Value = "'  -0.23%  "
$ws.Range("D8").Value = "'0.3092"
$ws.Range("E8").Value = "'  -2.76%  "
$ws.Range("D9").Value = "'26.51"
$ws.Range("E9").Value = "'  -5.02%  "
$ws.Range("D10").Value = "'0.06971"
$ws.Range("E10").Value = "'  -0.76%  "
$ws.Range("D11").Value = "'0.08077"
$ws.Range("E11").Value = "'  +0.41%  "
$ws.Range("D12").Value = "'0.7710"
$ws.Range("E12").Value = "'  -1.31%  "
$ws.Range("D13").Value = "'1.933.35"
$ws.Range("E13").Value = "'  -0.05%  "
$ws.Range("D14").Value = "'5.326"
$ws.Range("E14").Value = "'  -1.48%  "
$ws.Range("D15").Value = "'92.27"
$ws.Range("E15").Value = "'  -1.04%  "
$ws.Range("D16").Value = "'14.28"
$ws.Range("E16").Value = "'  -1.08%  "
$ws.Range("D17").Value = "'30.099.09"
$ws.Range("E17").Value = "'  -0.87%  "
$ws.Range("D18").Value = "'6.076"
$ws.Range("E18").Value = "'  +0.49%  "
$ws.Range("D19").Value = "'0.000007835"
$ws.Range("E19").Value = "'  -2.13%  "
$ws.Range("D20").Value = "'240.25"
$ws.Range("E20").Value = "'  -4.80%  "
$ws.Range("D21").Value = "'2.193.69"
$ws.Range("E21").Value = "'  +0.49%  "
$ws.Range("E22").Value = "'  -0.05%  "
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "'  -0.15%  "
$ws.Range("D24").Value = "'7.117"
$ws.Range("E24").Value = "'  +6.21%  "
$ws.Range("D25").Value = "'9.390"
$ws.Range("E25").Value = "'  -1.45%  "
$ws.Range("D26").Value = "'166.80"
$ws.Range("E26").Value = "'  +1.23%  "
$ws.Range("D27").Value = "'18.96"
$ws.Range("E27").Value = "'  -0.74%  "
$ws.Range("D28").Value = "'0.1275"
$ws.Range("E28").Value = "'  -1.89%  "
$ws.Range("D29").Value = "'2.045"
$ws.Range("E29").Value = "'  -8.02%  "
$ws.Range("E30").Value = "'  +1.25%  "
$ws.Range("D31").Value = "'1.350"
$ws.Range("E31").Value = "'  -1.83%  "
$ws.Range("D32").Value = "'4.339"
$ws.Range("D33").Value = "'4.085"
$ws.Range("E33").Value = "'  -1.26%  "
$ws.Range("D34").Value = "'1.310"
$ws.Range("E34").Value = "'  -2.14%  "
$ws.Range("D35").Value = "'0.05158"
$ws.Range("D36").Value = "'0.7498"
$ws.Range("E36").Value = "'  -0.93%  "
$ws.Range("D37").Value = "'2.721"
$ws.Range("E37").Value = "'  -2.56%  "
$ws.Range("D38").Value = "'0.01961"
$ws.Range("E38").Value = "'  +0.26%  "
$ws.Range("D39").Value = "'2.797"
$ws.Range("E39").Value = "'  -0.54%  "
$ws.Range("D40").Value = "'6.358"
$ws.Range("E40").Value = "'  -4.02%  "
$ws.Range("D41").Value = "'0.4504"
$ws.Range("E41").Value = "'  +0.07%  "
$ws.Range("D42").Value = "'74.54"
$ws.Range("D43").Value = "'1.981"
$ws.Range("E43").Value = "'  +0.15%  "
$ws.Range("E44").Value = "'  -0.10%  "
$ws.Range("D45").Value = "'0.8394"
$ws.Range("E45").Value = "'  +0.29%  "
$ws.Range("D46").Value = "'7.725"
$ws.Range("E46").Value = "'  +0.53%  "
$ws.Range("D47").Value = "'9.967"
$ws.Range("E47").Value = "'  -0.20%  "
$ws.Range("D48").Value = "'101.82"
$ws.Range("E48").Value = "'  +0.24%  "
$ws.Range("D49").Value = "'2.086.98"
$ws.Range("E49").Value = "'  +0.21%  "
$ws.Range("D50").Value = "'36.78"
$ws.Range("E50").Value = "'  -2.65%  "
$ws.Range("D51").Value = "'0.1183"
$ws.Range("E51").Value = "'  -4.63%  "
